$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7128567695617676
$ws.Range("E2").Value = 1725.274415146134
$ws.Range("F2").Value = 0.07714720742576983
$ws.Range("G2").Value = 0.05855277118619816
$ws.Range("H2").Value = 0.05301077648562695
$ws.Range("I2").Value = 0.04646098188441308
$ws.Range("J2").Value = 0.0448462350163345
$ws.Range("K2").Value = 0.04107815679574332
$ws.Range("L2").Value = 0.03992018699727849
$ws.Range("M2").Value = 0.03930612381589912
$ws.Range("N2").Value = 0.03732936222220197
$ws.Range("O2").Value = 0.03657275702526683
$ws.Range("P2").Value = 0.03613035294883613
$ws.Range("Q2").Value = 0.03534876826267415
$ws.Range("R2").Value = 0.03472844916699373
$ws.Range("S2").Value = 0.03447078037027117
$ws.Range("T2").Value = 0.03425956795038169
$ws.Range("U2").Value = 0.03398791345690835
$ws.Range("V2").Value = 0.03386674289939845
$ws.Range("W2").Value = 0.03365577351694862
$ws.Range("X2").Value = 0.03365577351694862
$ws.Range("Y2").Value = 0.03363108021727355
$ws.Range("C3").Value = 0.8094415664672852
$ws.Range("E3").Value = 1658.702412027558
$ws.Range("F3").Value = 0.07993603317253202
$ws.Range("G3").Value = 0.0600373871715624
$ws.Range("H3").Value = 0.05499061662251244
$ws.Range("I3").Value = 0.04898080373424343
$ws.Range("J3").Value = 0.04566936303030171
$ws.Range("K3").Value = 0.04175216857630489
$ws.Range("L3").Value = 0.03917546860501377
$ws.Range("M3").Value = 0.03799240527285993
$ws.Range("N3").Value = 0.03722204659710953
$ws.Range("O3").Value = 0.03614261932957002
$ws.Range("P3").Value = 0.03508027402640844
$ws.Range("Q3").Value = 0.03423759529042681
$ws.Range("R3").Value = 0.03402763590101018
$ws.Range("S3").Value = 0.03336594890418419
$ws.Range("T3").Value = 0.03317061747940166
$ws.Range("U3").Value = 0.03277150475498756
$ws.Range("V3").Value = 0.0326343308456293
$ws.Range("W3").Value = 0.03247985111898834
$ws.Range("X3").Value = 0.03247985111898834
$ws.Range("Y3").Value = 0.03233338035141438
$ws.Range("C4").Value = 0.7343530654907227
$ws.Range("E4").Value = 1708.415277644104
$ws.Range("F4").Value = 0.07606826902394109
$ws.Range("G4").Value = 0.06109773272867288
$ws.Range("H4").Value = 0.05391789002138476
$ws.Range("I4").Value = 0.04702797443863287
$ws.Range("J4").Value = 0.04559992039178696
$ws.Range("K4").Value = 0.04253952008099317
$ws.Range("L4").Value = 0.04002223933341751
$ws.Range("M4").Value = 0.03802215979607194
$ws.Range("N4").Value = 0.03712689967891258
$ws.Range("O4").Value = 0.03615670364545857
$ws.Range("P4").Value = 0.03584115536853046
$ws.Range("Q4").Value = 0.03545351742542135
$ws.Range("R4").Value = 0.0347461991067895
$ws.Range("S4").Value = 0.03434042495901449
$ws.Range("T4").Value = 0.03417704843906031
$ws.Range("U4").Value = 0.03397376352729545
$ws.Range("V4").Value = 0.03365718583948557
$ws.Range("W4").Value = 0.03357436909643202
$ws.Range("X4").Value = 0.03345316569715985
$ws.Range("Y4").Value = 0.03330244205933926
$ws.Range("C5").Value = 0.7343335151672363
$ws.Range("E5").Value = 1693.299554676385
$ws.Range("F5").Value = 0.07725229382936318
$ws.Range("G5").Value = 0.06224267556445609
$ws.Range("H5").Value = 0.05337554943318661
$ws.Range("I5").Value = 0.04813063474834807
$ws.Range("J5").Value = 0.04411049026703234
$ws.Range("K5").Value = 0.04157367011588168
$ws.Range("L5").Value = 0.03937308822630312
$ws.Range("M5").Value = 0.03793309016194655
$ws.Range("N5").Value = 0.03669218826340355
$ws.Range("O5").Value = 0.03617825122849138
$ws.Range("P5").Value = 0.03534500346699983
$ws.Range("Q5").Value = 0.03486109752489246
$ws.Range("R5").Value = 0.03437229913363837
$ws.Range("S5").Value = 0.03379274364108814
$ws.Range("T5").Value = 0.0337254397876894
$ws.Range("U5").Value = 0.03346200741746601
$ws.Range("V5").Value = 0.03335871451347492
$ws.Range("W5").Value = 0.03318696676778664
$ws.Range("X5").Value = 0.03300778859018293
$ws.Range("Y5").Value = 0.03300778859018293
$ws.Range("C6").Value = 0.7187759876251221
$ws.Range("E6").Value = 1704.847021420177
$ws.Range("F6").Value = 0.07794159007063568
$ws.Range("G6").Value = 0.06329787779973907
$ws.Range("H6").Value = 0.05550485021193578
$ws.Range("I6").Value = 0.0478588651650704
$ws.Range("J6").Value = 0.045242900820936
$ws.Range("K6").Value = 0.04234436538865143
$ws.Range("L6").Value = 0.03958260295335758
$ws.Range("M6").Value = 0.03872685471497382
$ws.Range("N6").Value = 0.03754319307850806
$ws.Range("O6").Value = 0.03675848115232794
$ws.Range("P6").Value = 0.0356983340963338
$ws.Range("Q6").Value = 0.0350778915528709
$ws.Range("R6").Value = 0.03462047877221776
$ws.Range("S6").Value = 0.03432666219765652
$ws.Range("T6").Value = 0.03395924391796412
$ws.Range("U6").Value = 0.03376242267105163
$ws.Range("V6").Value = 0.03359472030223437
$ws.Range("W6").Value = 0.03342741820246359
$ws.Range("X6").Value = 0.03331135722652187
$ws.Range("Y6").Value = 0.03323288540780071
$ws.Range("C7").Value = 0.7030963897705078
$ws.Range("E7").Value = 1728.472235158273
$ws.Range("F7").Value = 0.07940480848271377
$ws.Range("G7").Value = 0.06451254469356596
$ws.Range("H7").Value = 0.05549503521488579
$ws.Range("I7").Value = 0.05103093036727432
$ws.Range("J7").Value = 0.04610768096234179
$ws.Range("K7").Value = 0.04299541782844442
$ws.Range("L7").Value = 0.04244156212311004
$ws.Range("M7").Value = 0.04055013089345384
$ws.Range("N7").Value = 0.03829293635685287
$ws.Range("O7").Value = 0.03736744981178072
$ws.Range("P7").Value = 0.03621178293317352
$ws.Range("Q7").Value = 0.0356466010330479
$ws.Range("R7").Value = 0.03528309583004246
$ws.Range("S7").Value = 0.03486234289254527
$ws.Range("T7").Value = 0.03442278249858074
$ws.Range("U7").Value = 0.03413150379194606
$ws.Range("V7").Value = 0.03402147642345688
$ws.Range("W7").Value = 0.03391831949584324
$ws.Range("X7").Value = 0.03382486767686352
$ws.Range("Y7").Value = 0.03369341589002481
$ws.Range("C8").Value = 0.6875088214874268
$ws.Range("E8").Value = 1680.038990442041
$ws.Range("F8").Value = 0.07830246685610548
$ws.Range("G8").Value = 0.06283869090745071
$ws.Range("H8").Value = 0.0534260486328658
$ws.Range("I8").Value = 0.04920615319165435
$ws.Range("J8").Value = 0.04484476874954219
$ws.Range("K8").Value = 0.04229511819475675
$ws.Range("L8").Value = 0.03861381353294233
$ws.Range("M8").Value = 0.03691677307593486
$ws.Range("N8").Value = 0.03691677307593486
$ws.Range("O8").Value = 0.03590062314931426
$ws.Range("P8").Value = 0.03485282887654798
$ws.Range("Q8").Value = 0.03443334513409598
$ws.Range("R8").Value = 0.03404450040945926
$ws.Range("S8").Value = 0.03360131453048135
$ws.Range("T8").Value = 0.03343387376591433
$ws.Range("U8").Value = 0.03328303567510473
$ws.Range("V8").Value = 0.03304194428298216
$ws.Range("W8").Value = 0.03285499395544138
$ws.Range("X8").Value = 0.03280799244900118
$ws.Range("Y8").Value = 0.03274929805929903
$ws.Range("C9").Value = 0.703129768371582
$ws.Range("E9").Value = 1775.72829100966
$ws.Range("F9").Value = 0.07880690255613798
$ws.Range("G9").Value = 0.06506375271980157
$ws.Range("H9").Value = 0.05696212099949879
$ws.Range("I9").Value = 0.05207712022951443
$ws.Range("J9").Value = 0.04578365913553175
$ws.Range("K9").Value = 0.04387744728167087
$ws.Range("L9").Value = 0.042518893195849
$ws.Range("M9").Value = 0.0400060634507805
$ws.Range("N9").Value = 0.03861862295610691
$ws.Range("O9").Value = 0.03816059024430382
$ws.Range("P9").Value = 0.03711314373472037
$ws.Range("Q9").Value = 0.03673187201355669
$ws.Range("R9").Value = 0.03641237382955771
$ws.Range("S9").Value = 0.03602062815042889
$ws.Range("T9").Value = 0.03570080141135774
$ws.Range("U9").Value = 0.03528960678878811
$ws.Range("V9").Value = 0.03505647576158989
$ws.Range("W9").Value = 0.03484362310194385
$ws.Range("X9").Value = 0.03469003143531075
$ws.Range("Y9").Value = 0.03461458656938908
$ws.Range("C10").Value = 0.7031209468841553
$ws.Range("E10").Value = 1691.703501315402
$ws.Range("F10").Value = 0.07842689643623364
$ws.Range("G10").Value = 0.0621331276193114
$ws.Range("H10").Value = 0.05310217652751625
$ws.Range("I10").Value = 0.0483177164363271
$ws.Range("J10").Value = 0.04401997426853214
$ws.Range("K10").Value = 0.04136099050582896
$ws.Range("L10").Value = 0.03995021014167054
$ws.Range("M10").Value = 0.03785608032918315
$ws.Range("N10").Value = 0.03697074912238177
$ws.Range("O10").Value = 0.03568825941113152
$ws.Range("P10").Value = 0.03498220847983304
$ws.Range("Q10").Value = 0.03469506955567896
$ws.Range("R10").Value = 0.03433306585524827
$ws.Range("S10").Value = 0.034081061141456
$ws.Range("T10").Value = 0.03379302389467552
$ws.Range("U10").Value = 0.03356591448854521
$ws.Range("V10").Value = 0.03334125803678796
$ws.Range("W10").Value = 0.03311105877208755
$ws.Range("X10").Value = 0.03303681610980413
$ws.Range("Y10").Value = 0.03297667643889671
$ws.Range("C11").Value = 0.7031245231628418
$ws.Range("E11").Value = 1653.403176580196
$ws.Range("F11").Value = 0.07871652169297166
$ws.Range("G11").Value = 0.06354852415053422
$ws.Range("H11").Value = 0.05444668228980518
$ws.Range("I11").Value = 0.04881723221947781
$ws.Range("J11").Value = 0.04491349126484675
$ws.Range("K11").Value = 0.04254096343098179
$ws.Range("L11").Value = 0.03962420679707412
$ws.Range("M11").Value = 0.03842751199838352
$ws.Range("N11").Value = 0.03703793615402241
$ws.Range("O11").Value = 0.03595121227445278
$ws.Range("P11").Value = 0.03471399195533352
$ws.Range("Q11").Value = 0.0344641940596746
$ws.Range("R11").Value = 0.03370611772310028
$ws.Range("S11").Value = 0.03334835980723699
$ws.Range("T11").Value = 0.03315031460488575
$ws.Range("U11").Value = 0.03272032019017521
$ws.Range("V11").Value = 0.03248284800602291
$ws.Range("W11").Value = 0.03241957335844484
$ws.Range("X11").Value = 0.03231173022168993
$ws.Range("Y11").Value = 0.03223008141481862
